# hardware, ordering: remove knobs from expansion panel bom
#
# The "Chassis Pod" section of the BOM included two knob line items
# (Item # 208 and 209, manufactured by Davies Molding, LLC) that are no
# longer used on the expansion panel. Remove those two rows entirely;
# everything below shifts up to close the gap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows 20 and 21 are the two knob entries under the "Chassis Pod" header.
$ws.Rows.Item(20).Resize(2).Delete()

# Leave the sheet scrolled/selected near the top of the data, similar to
# where the editor's cursor ended up after removing the rows.
$ws.Range("A20:XFD20").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
